$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").Value = 19999996249
$ws.Range("F5").Value = 30193317188

$ws.Application.ActiveWindow.ScrollColumn = 6
$ws.Range("M21").Select()

$wb.Application.Calculate()
